{"js": "// Insert the GitHub repo URL right after the existing\n// \"My code can be found at \" sentence, matching the author's commit\n// \"include url to github repo\".\nconst searchText = \"My code can be found at \";\nconst url = \"https://github.com/mdamiani610/SlowLifeGUI\";\n\nconst results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (!results.items.length) {\n  throw new Error(`Could not find target text \"${searchText}\" in the document body.`);\n}\n\n// There is a single occurrence in this document; take the first match and\n// insert the URL immediately after it (same run formatting carries over).\nconst target = results.items[0];\ntarget.insertText(url, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Insert the GitHub repo URL right after the existing\n# \"My code can be found at \" sentence, matching the author's commit\n# \"include url to github repo\".\n$d = $word.ActiveDocument\n\n$searchText = \"My code can be found at \"\n$url = \"https://github.com/mdamiani610/SlowLifeGUI\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = $searchText\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.MatchWildcards = $false\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n\n$found = $rng.Find.Execute()\n\nif ($found) {\n    # Collapse the found range to its end point, then insert the URL there\n    # (inherits the formatting of the text it's inserted next to).\n    $rng.Collapse(0)\n    $rng.InsertAfter($url)\n} else {\n    throw \"Could not find target text '$searchText' in the document.\"\n}\n"}
